$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.794.65'
$ws.Range('E2').Value = '  +2.03%  '
$ws.Range('D3').Value = '3.271.17'
$ws.Range('E3').Value = '  -0.63%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '567.94'
$ws.Range('E5').Value = '  -2.00%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '175.82'
$ws.Range('E6').Value = '  -3.89%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  +2.08%  '
$ws.Range('D9').Value = '3.261.96'
$ws.Range('E9').Value = '  -0.74%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.174'
$ws.Range('E11').Value = '  +0.33%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '45.53'
$ws.Range('E12').Value = '  -1.92%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000269'
$ws.Range('E13').Value = '  +2.12%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '691.85'
$ws.Range('E14').Value = '  +9.39%  '
$ws.Range('D15').Value = '3.796.25'
$ws.Range('E15').Value = '  -0.50%  '
$ws.Range('E16').Value = '  -1.31%  '
$ws.Range('D17').Value = '66.851.85'
$ws.Range('E17').Value = '  +1.93%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.119'
$ws.Range('E18').Value = '  +1.02%  '
$ws.Range('D19').Value = '3.278.22'
$ws.Range('E19').Value = '  -0.30%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '17.31'
$ws.Range('E20').Value = '  -2.04%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '10.74'
$ws.Range('E21').Value = '  -1.78%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.886'
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('E23').Value = '  -5.67%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '97.81'
$ws.Range('E25').Value = '  -2.66%  '
$ws.Range('E26').Value = '  -2.37%  '
$ws.Range('E27').Value = '  -1.27%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.29'
$ws.Range('E28').Value = '  -0.50%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '32.73'
$ws.Range('E29').Value = '  +6.66%  '
$ws.Range('E30').Value = '  +0.80%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.74'
$ws.Range('E31').Value = '  +3.86%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '570.13'
$ws.Range('E32').Value = '  -0.84%  '
$ws.Range('D33').Value = '3.870.28'
$ws.Range('E33').Value = '  +0.61%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '10.79'
$ws.Range('E34').Value = '  -0.57%  '
$ws.Range('E35').Value = '  -0.95%  '
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('B37').Value = 'dogwifhat'
$ws.Range('C37').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.32'
$ws.Range('E37').Value = '  -9.10%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '55.26'
$ws.Range('E38').Value = '  -0.38%  '
$ws.Range('E39').Value = '  +2.12%  '
$ws.Range('E40').Value = '  +0.62%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.35'
$ws.Range('E41').Value = '  -1.70%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '31.68'
$ws.Range('E42').Value = '  -2.10%  '
$ws.Range('D43').Value = '0.0₃0670'
$ws.Range('E43').Value = '  -1.69%  '
$ws.Range('E44').Value = '  -2.44%  '
$ws.Range('E45').Value = '  -1.34%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0406'
$ws.Range('E46').Value = '  +0.45%  '
$ws.Range('E47').Value = '  +0.64%  '
$ws.Range('E48').Value = '  +0.19%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.53'
$ws.Range('E49').Value = '  +1.12%  '
$ws.Range('E50').Value = '  +8.08%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '129.47'
$ws.Range('E51').Value = '  +0.22%  '
